$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")

# --- Fill in the new Main-sheet rows (CE002..CE008) -------------------------
# Shared-string insertion order matters for exact OOXML parity: write all of
# column B (descriptions) first, then column A (part numbers), then column C
# (branch), matching the order the strings appear in the target file.

$main.Range("B5").Value  = "External Storage"
$main.Range("B6").Value  = "Solenoid Valve"
$main.Range("B7").Value  = "Lithium ion Battery"
$main.Range("B8").Value  = "Pump"
$main.Range("B9").Value  = "Spark igniter"
$main.Range("B10").Value = "Thermocouple"
$main.Range("B11").Value = "Pressure Transducer"

$main.Range("A5").Value  = "CE002"
$main.Range("A6").Value  = "CE003"
$main.Range("A7").Value  = "CE004"
$main.Range("A8").Value  = "CE005"
$main.Range("A9").Value  = "CE006"
$main.Range("A10").Value = "CE007"
$main.Range("A11").Value = "CE008"

$main.Range("C5").Value  = "Avionics"
$main.Range("C6").Value  = "Avionics"
$main.Range("C7").Value  = "Avionics"
$main.Range("C8").Value  = "Avionics"
$main.Range("C9").Value  = "Avionics"
$main.Range("C10").Value = "Avionics"
$main.Range("C11").Value = "Avionics"

# D5 carries the Hyperlink visual style only (no value/link), same as D4's style.
$main.Range("D5").Style = "Hyperlink"

# --- Add the new per-component worksheets (CE002..CE007) --------------------
$names = @("CE002", "CE003", "CE004", "CE005", "CE006", "CE007")

foreach ($name in $names) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $name
    $ws.Range("A1").Value = "Part"
    $ws.Range("B1").Value = $name
    $ws.Range("A1:B1").Font.Bold = $true
    $ws.Range("A1:B1").Font.Size = 16
    $ws.Range("B1").Select()
}

# Last new sheet (CE007) keeps a stray selection at L6 in the target file.
$ce007 = $wb.Worksheets.Item("CE007")
$ce007.Range("L6").Select()

# --- Restore view-state on the pre-existing sheets ---------------------------
$ce001 = $wb.Worksheets.Item("CE001")
$ce001.Range("A1").Select()

$main.Activate()
$main.Range("D5").Select()
